$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column B ("ASIN" etc. all shift one column right)
$ws.Columns.Item(2).Insert()

# New column B header + values (Week_Start_Date) - force text so the
# ISO date strings are not auto-converted to Excel date serials
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"
$ws.Range("B2:B17").NumberFormat = "@"

$dates = @(
    "2025-01-05", "2025-01-12", "2025-01-19", "2025-01-26",
    "2025-02-02", "2025-02-09", "2025-02-16", "2025-02-23",
    "2025-03-02", "2025-03-09", "2025-03-16", "2025-03-23",
    "2025-03-30", "2025-04-06", "2025-04-13", "2025-04-20"
)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $dates[$i]
}

# Column A week labels drop the leading zero (W01 -> W1 ... W09 -> W9);
# W10..W16 already match and are left untouched
$weeks = @("W1", "W2", "W3", "W4", "W5", "W6", "W7", "W8", "W9")
for ($i = 0; $i -lt $weeks.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $weeks[$i]
}

# is_holiday_week (now column J after the insert) becomes a real boolean
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 10).Value = $false
}
